# update subject fee set
$p = $ppt.ActivePresentation

# --- Slide with "基础被试费：	+70元" -> "基础被试费：	+60元" ---
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(1)
$tr12 = $sh12.TextFrame.TextRange
$tr12.Characters(17, 3).Text = "+60"

# --- Slide with the green/red "色" explanation runs ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(3)
$tr3 = $sh3.TextFrame.TextRange
# "色" + "表示得到" + "，"  ->  single run "色表示得到，"
$tr3.Characters(41, 6).Text = "色表示得到，"
# "色" + "表示失去"  ->  single run "色表示失去"
$tr3.Characters(48, 5).Text = "色表示失去"

# --- Slide with "...抽出一个实现，..." -> "...抽出两个实现，..." ---
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$tr5 = $sh5.TextFrame.TextRange
$tr5.Characters(71, 3).Text = "出两个"
